$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.301.25'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.08%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.704.59'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.03%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9975'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.31'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.60%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9976'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.21%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4872'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.23%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2575'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06153'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.20%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.714.20'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.47%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06939'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.69%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.45'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.96%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.5956'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.23%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.443'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.20%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.23'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.25%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9972'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.20%  '

$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9982'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.12%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.227.00'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.40%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007037'
$ws.Range('D19').Style = 'Normal'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.41%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.932.81'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.35%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.372'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.91%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.362'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.41%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.001'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.84%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.77%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.14'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.15%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.389'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.15%  '

$ws.Range('B28').Value = 'BitcoinCash'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '105.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.49%  '

$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.720'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.64%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.863'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.72%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07913'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.06%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.582'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.99%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04434'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.28%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.601'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.40%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9870'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.31%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6130'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.94%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9505'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.48%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.984'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.53%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.351'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.15%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9955'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.95%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01470'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.42%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.28'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.56%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.362'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.94%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3781'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.96%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.766'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.18%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1143'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.86%  '

$ws.Range('E47').Value = '  -1.01%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.34'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.96%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.704'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.72%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '50.95'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.64%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9990'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.34%  '
